$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header B1: "locacalizacion" -> "localizacion"
$ws.Range("B1").Value = "localizacion"

# Clear the content of D2 (previously "ID4"), keep its style
$ws.Range("D2").Value = ""

# Update selection to active cell D2
$ws.Range("D2").Select()
